$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Step 1: Insert the new "2022-Q3" fund-holdings detail sheet.
#   Easiest way to get an exact layout/style match (bold bordered header
#   row, bordered index column) is to clone the existing "2022-Q2" sheet
#   (same columns/headers/styles) and place the clone right in front of
#   it, then rename the clone and overwrite its data.
# -------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($src, $null)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# "2022-Q2" had 21 data rows (1 header + 20 funds); "2022-Q3" only has
# 9 (1 header + 8 funds), so drop the extra cloned rows.
$newSheet.Rows("10:21").Delete()

# Columns B-G store fund code / name / scale / position% / ratio% / value
# as TEXT (e.g. "001475", "164.58") even though several look numeric -
# force Text format before writing so Excel doesn't silently coerce them
# to numbers (which would also strip the leading zeros from fund codes).
$newSheet.Range("B2:G9").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "001475"
$newSheet.Cells.Item(2, 3).Value = "易方达国防军工混合A"
$newSheet.Cells.Item(2, 4).Value = "164.58"
$newSheet.Cells.Item(2, 5).Value = "93.30"
$newSheet.Cells.Item(2, 6).Value = "3.80"
$newSheet.Cells.Item(2, 7).Value = "6.2540"
$newSheet.Cells.Item(2, 8).Value = 10

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "001838"
$newSheet.Cells.Item(3, 3).Value = "国投瑞银国家安全灵活配置混合"
$newSheet.Cells.Item(3, 4).Value = "27.40"
$newSheet.Cells.Item(3, 5).Value = "94.42"
$newSheet.Cells.Item(3, 6).Value = "8.89"
$newSheet.Cells.Item(3, 7).Value = "2.4359"
$newSheet.Cells.Item(3, 8).Value = 4

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "005774"
$newSheet.Cells.Item(4, 3).Value = "华夏产业升级混合A"
$newSheet.Cells.Item(4, 4).Value = "24.29"
$newSheet.Cells.Item(4, 5).Value = "93.85"
$newSheet.Cells.Item(4, 6).Value = "5.72"
$newSheet.Cells.Item(4, 7).Value = "1.3894"
$newSheet.Cells.Item(4, 8).Value = 9

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "015059"
$newSheet.Cells.Item(5, 3).Value = "华夏产业升级混合C"
$newSheet.Cells.Item(5, 4).Value = "8.92"
$newSheet.Cells.Item(5, 5).Value = "93.85"
$newSheet.Cells.Item(5, 6).Value = "5.72"
$newSheet.Cells.Item(5, 7).Value = "0.5102"
$newSheet.Cells.Item(5, 8).Value = 9

$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).Value = "015945"
$newSheet.Cells.Item(6, 3).Value = "易方达国防军工混合C"
$newSheet.Cells.Item(6, 4).Value = "10.80"
$newSheet.Cells.Item(6, 5).Value = "93.30"
$newSheet.Cells.Item(6, 6).Value = "3.80"
$newSheet.Cells.Item(6, 7).Value = "0.4104"
$newSheet.Cells.Item(6, 8).Value = 10

$newSheet.Cells.Item(7, 1).Value = 5
$newSheet.Cells.Item(7, 2).Value = "012239"
$newSheet.Cells.Item(7, 3).Value = "惠升优势企业一年持有期灵活配置混合"
$newSheet.Cells.Item(7, 4).Value = "10.70"
$newSheet.Cells.Item(7, 5).Value = "79.97"
$newSheet.Cells.Item(7, 6).Value = "3.42"
$newSheet.Cells.Item(7, 7).Value = "0.3659"
$newSheet.Cells.Item(7, 8).Value = 6

$newSheet.Cells.Item(8, 1).Value = 6
$newSheet.Cells.Item(8, 2).Value = "012568"
$newSheet.Cells.Item(8, 3).Value = "天弘高端制造混合A"
$newSheet.Cells.Item(8, 4).Value = "5.98"
$newSheet.Cells.Item(8, 5).Value = "92.42"
$newSheet.Cells.Item(8, 6).Value = "4.61"
$newSheet.Cells.Item(8, 7).Value = "0.2757"
$newSheet.Cells.Item(8, 8).Value = 7

$newSheet.Cells.Item(9, 1).Value = 7
$newSheet.Cells.Item(9, 2).Value = "012569"
$newSheet.Cells.Item(9, 3).Value = "天弘高端制造混合C"
$newSheet.Cells.Item(9, 4).Value = "0.97"
$newSheet.Cells.Item(9, 5).Value = "92.42"
$newSheet.Cells.Item(9, 6).Value = "4.61"
$newSheet.Cells.Item(9, 7).Value = "0.0447"
$newSheet.Cells.Item(9, 8).Value = 7

# -------------------------------------------------------------------------
# Step 2: Update the "总计" (totals) sheet - insert a new row for 2022-Q3
# right under the header, pushing the existing quarters down one row,
# and bump every row's running index in column A by 1.
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows(2).Insert()

# Row-insert copies formatting from the row above (the bold header),
# which the data rows should not have - clear B2:D2 back to plain, then
# restore column A's "index" cell style by copying it from A3 (which
# already carries the correct bordered/centered style).
$ws.Range("B2:D2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 11.69

for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "Added 2022-Q3 sheet and updated 总计 summary."
